$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.943.16"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "1.673.23"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("E6").Value = "  +1.74%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E8").Value = "  -0.29%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.25"
$ws.Range("E10").Value = "  +0.49%  "

$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").Value = "1.909.16"
$ws.Range("E12").Value = "  +0.82%  "

$ws.Range("D13").Value = "1.678.86"
$ws.Range("E13").Value = "  +1.12%  "

$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.58"
$ws.Range("E16").Value = "  -0.21%  "

$ws.Range("D17").Value = "26.938.88"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.09"
$ws.Range("E18").Value = "  +4.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "234.97"
$ws.Range("E19").Value = "  -0.77%  "

$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("E23").Value = "  -1.30%  "

$ws.Range("E24").Value = "  -2.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.43"
$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.19"
$ws.Range("E26").Value = "  +0.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.00"
$ws.Range("E27").Value = "  +0.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("E28").Value = "  -1.25%  "

$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").Value = "1.476.90"
$ws.Range("E33").Value = "  -4.66%  "

$ws.Range("E34").Value = "  +1.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.67"
$ws.Range("E35").Value = "  +2.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("E37").Value = "  +1.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.894"
$ws.Range("E38").Value = "  -1.33%  "

$ws.Range("E39").Value = "  +0.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.04"
$ws.Range("E40").Value = "  +7.53%  "

$ws.Range("E41").Value = "  -3.58%  "

$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("E43").Value = "  +3.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.39"
$ws.Range("E44").Value = "  +0.99%  "

$ws.Range("D45").Value = "1.814.54"
$ws.Range("E45").Value = "  +0.59%  "

$ws.Range("E46").Value = "  -0.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.48"
$ws.Range("E47").Value = "  -0.18%  "

$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("E49").Value = "  +1.52%  "

# Row 50: EnergySwap -> Cronos
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0508"
$ws.Range("E50").Value = "  +0.26%  "

# Row 51: Cronos -> EnergySwap
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.72"
$ws.Range("E51").Value = "  +0.22%  "

Write-Output "Applied all changes"